$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Test")
$ws.Activate()

# Row 5 ("SES-02") gets updated to a new "SES-TC-02 / SES-TS-02" test case
# (verify session termination after logout), replacing the previous
# "LOGOUT-TC-02 / LOGOUT-TS-02" content that used to live there.
$ws.Range("E5").Value = "SES-TC-02"
$ws.Range("G5").Value = "SES-TS-02"
$ws.Range("D5").Value = "Verify session is terminated after logout"
$ws.Range("F5").Value = "Session Module"
$ws.Range("H5").Value = "Verify session is destroyed after logout"
$ws.Range("I5").Value = "User already logged in"
$ws.Range("J5").Value = "1. Click Logout `n2. Try to refresh page `n3. Try to access restricted page"
$ws.Range("K5").Value = "-"
$ws.Range("L5").Value = "Session is invalidated and user is redirected to login page"

# Row height recalculated (shorter wrapped text now fits in fewer lines)
$ws.Rows.Item(5).RowHeight = 47.25

# Columns H and I now share the same (wider) best-fit width
$ws.Columns.Item(8).ColumnWidth = 34.33

# Update the active selection left after the edit
$ws.Range("J6").Select() | Out-Null
